# Generate Report for Handback
#
# The c7dafe61-... file has now been handed back (in sync with en-US), so the
# status/report rows rotate: c7dafe61 moves to the top (row 2) with fresh
# handback data, ffff25e0... moves down to row 3, and ffffffbbafeb98... moves
# down to row 4 (taking the "duplicate content" flag with it).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
$ws1.Range("B2").Value = "e2e\c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
$ws1.Range("G2").Value = "2016-09-01 09:22:19"

$ws1.Range("A3").Value = "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
$ws1.Range("B3").Value = "e2e\ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"

$ws1.Range("A4").Value = "ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md"
$ws1.Range("B4").Value = "e2e\ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md"
$ws1.Range("E4").Value = "Handed back: in sync with en-US"
$ws1.Range("F4").Value = "Handed back: in sync with en-US"
$ws1.Range("G4").Value = "2016-09-01 09:20:15"

# Hyperlinks follow the same rows (B2/B3/B4) but now point at the file that
# now occupies that row, so rebuild them in the new order.
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7175abd06daa152eda906890f5c6df7854560197/e2e/c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md", "", "", "e2e\c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eae6ec605ac9b1f5bba0fc696da76dbac6e65288/e2e/ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md", "", "", "e2e\ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md") | Out-Null
$ws1.Hyperlinks.Add($ws1.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7175abd06daa152eda906890f5c6df7854560197/e2e/ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md", "", "", "e2e\ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
$ws2.Range("G2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.1d934907a74ac1423164f5eb0eb4fb60ad0e471f.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-09-01 09:22:14"
$ws2.Range("I2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
$ws2.Range("J2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.1d934907a74ac1423164f5eb0eb4fb60ad0e471f.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-09-01 09:22:40"

$ws2.Range("A3").Value = "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
$ws2.Range("F3").Value = "False"

$ws2.Range("A4").Value = "ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md"
$ws2.Range("C4").Value = "Handed back: in sync with en-US"
$ws2.Range("F4").Value = "True"
$ws2.Range("G4").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.zh-cn.xlf"
$ws2.Range("H4").Value = "2016-09-01 09:19:58"
$ws2.Range("I4").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
$ws2.Range("J4").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.zh-cn.xlf"
$ws2.Range("K4").Value = "2016-09-01 09:20:36"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7175abd06daa152eda906890f5c6df7854560197/e2e/c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md", "", "", "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0a72c1baf1d40f6ddfe6878c79b79b43da756dca/e2e/836d1c83-52d1-4579-9d23-7f3bdff4659d.md", "", "", "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eae6ec605ac9b1f5bba0fc696da76dbac6e65288/e2e/ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md", "", "", "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0a72c1baf1d40f6ddfe6878c79b79b43da756dca/e2e/836d1c83-52d1-4579-9d23-7f3bdff4659d.md", "", "", "836d1c83-52d1-4579-9d23-7f3bdff4659d.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7175abd06daa152eda906890f5c6df7854560197/e2e/ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md", "", "", "ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md") | Out-Null
$ws2.Hyperlinks.Add($ws2.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/0a72c1baf1d40f6ddfe6878c79b79b43da756dca/e2e/836d1c83-52d1-4579-9d23-7f3bdff4659d.md", "", "", "836d1c83-52d1-4579-9d23-7f3bdff4659d.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
$ws3.Range("G2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.1d934907a74ac1423164f5eb0eb4fb60ad0e471f.de-de.xlf"
$ws3.Range("H2").Value = "2016-09-01 09:22:19"
$ws3.Range("I2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md"
$ws3.Range("J2").Value = "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.1d934907a74ac1423164f5eb0eb4fb60ad0e471f.de-de.xlf"
$ws3.Range("K2").Value = "2016-09-01 09:22:49"

$ws3.Range("A3").Value = "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md"
$ws3.Range("F3").Value = "False"

$ws3.Range("A4").Value = "ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md"
$ws3.Range("C4").Value = "Handed back: in sync with en-US"
$ws3.Range("F4").Value = "True"
$ws3.Range("G4").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.de-de.xlf"
$ws3.Range("H4").Value = "2016-09-01 09:20:15"
$ws3.Range("I4").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.md"
$ws3.Range("J4").Value = "836d1c83-52d1-4579-9d23-7f3bdff4659d.03b0573f4c3ad94c7c118a0f76853d51b4646da2.de-de.xlf"
$ws3.Range("K4").Value = "2016-09-01 09:20:43"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7175abd06daa152eda906890f5c6df7854560197/e2e/c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md", "", "", "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6da3816f755bbc0b843c2d8a6cb88b2d91903c8c/e2e/836d1c83-52d1-4579-9d23-7f3bdff4659d.md", "", "", "c7dafe61-f43a-4dfb-aa31-5ae2eae150d1.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eae6ec605ac9b1f5bba0fc696da76dbac6e65288/e2e/ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md", "", "", "ffff25e0add6-8617-465c-9af1-a53a3d57b9a0.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6da3816f755bbc0b843c2d8a6cb88b2d91903c8c/e2e/836d1c83-52d1-4579-9d23-7f3bdff4659d.md", "", "", "836d1c83-52d1-4579-9d23-7f3bdff4659d.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7175abd06daa152eda906890f5c6df7854560197/e2e/ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md", "", "", "ffffffbbafeb98-b458-4f5c-9134-b345a406d458.md") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6da3816f755bbc0b843c2d8a6cb88b2d91903c8c/e2e/836d1c83-52d1-4579-9d23-7f3bdff4659d.md", "", "", "836d1c83-52d1-4579-9d23-7f3bdff4659d.md") | Out-Null
